{"js": "// Prepend \"Design: \" to the start of each Q&A \"answer\" bullet paragraph\n// in the feedback table (the six ListBullet paragraphs that previously\n// had no such prefix).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// Exact original texts (start of paragraph) that must receive the\n// \"Design: \" prefix. Matching is done on a prefix of the paragraph's\n// text so paragraphs containing a <w:br/> (multiple runs) still match\n// on their first run's content.\nconst targets = [\n  \"Mostly\",\n  \"Nothing special to mention\",\n  \"91xx Went well, some minor budget challenges\",\n  \"Mostly ok.\",\n  \"Internal communication ok. External communication with suppliers mostly ok.\",\n  \"Some things went to correct direction but regarding TK the opposite way.\",\n];\n\nconst prefix = \"Design: \";\n\nfor (const para of paragraphs.items) {\n  const text = para.text;\n  for (const target of targets) {\n    const alreadyPrefixed = text.indexOf(prefix + target) === 0;\n    const isTarget = text.indexOf(target) === 0;\n    if (isTarget && !alreadyPrefixed) {\n      para.insertText(prefix, Word.InsertLocation.start);\n      break;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Prepend \"Design: \" to the start of each Q&A \"answer\" bullet paragraph\n# in the feedback table (the six ListBullet paragraphs that previously\n# had no such prefix).\n$d = $word.ActiveDocument\n\n$targets = @(\n    \"Mostly\",\n    \"Nothing special to mention\",\n    \"91xx Went well, some minor budget challenges\",\n    \"Mostly ok.\",\n    \"Internal communication ok. External communication with suppliers mostly ok.\",\n    \"Some things went to correct direction but regarding TK the opposite way.\"\n)\n\n$prefix = \"Design: \"\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    # Drop the trailing paragraph-mark / cell-mark control character(s)\n    # before comparing, so the match is against the visible text only.\n    $trimmed = $text.TrimEnd([char]13, [char]7)\n\n    foreach ($target in $targets) {\n        if ($trimmed.StartsWith($prefix + $target)) {\n            break\n        }\n        if ($trimmed.StartsWith($target)) {\n            $p.Range.InsertBefore($prefix)\n            break\n        }\n    }\n}\n"}
